$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 97.51272252350722
$ws.Range("C3").Value = 97.69870314573896
$ws.Range("C4").Value = 98.20421998602659
$ws.Range("C5").Value = 99.21022586325735
$ws.Range("C6").Value = 99.58756543339153
$ws.Range("C7").Value = 99.71841607293898
$ws.Range("C8").Value = 99.22934341072109
$ws.Range("C9").Value = 99.18395602243879
$ws.Range("C10").Value = 98.95500188459366
$ws.Range("C11").Value = 99.48113949271229
$ws.Range("C12").Value = 99.83825399215294
$ws.Range("C13").Value = 100.0050260107236
$ws.Range("C14").Value = 99.23934757934764
$ws.Range("C15").Value = 99.3796086926158
$ws.Range("C16").Value = 99.31481018560939
$ws.Range("C17").Value = 99.77996896572293
$ws.Range("C18").Value = 100.0077876080712
$ws.Range("C19").Value = 100.5948674333569
$ws.Range("C20").Value = 99.89564115076296
$ws.Range("C21").Value = 99.9294829530583
$ws.Range("C22").Value = 99.80972244276765
$ws.Range("C23").Value = 100.2130885795517
$ws.Range("C24").Value = 100.3388879102719
$ws.Range("C25").Value = 100.2683904079147
$ws.Range("C27").Value = 99.29312052197895
$ws.Range("C28").Value = 99.99921286471898
$ws.Range("C29").Value = 100.867268066879
$ws.Range("C30").Value = 101.2169134424177
$ws.Range("C31").Value = 101.5712489272691
$ws.Range("C32").Value = 100.7300852772187
$ws.Range("C33").Value = 101.0688688580645
$ws.Range("C34").Value = 101.4193216989077
$ws.Range("C35").Value = 102.8090123587009
$ws.Range("C36").Value = 103.043258537902
$ws.Range("C37").Value = 103.9630198827306
$ws.Range("C38").Value = 103.9540581538035
$ws.Range("C39").Value = 104.5673225817185
$ws.Range("C40").Value = 107.1860661179801
$ws.Range("C41").Value = 107.0401595331048
$ws.Range("C42").Value = 107.9805007201712
$ws.Range("C43").Value = 110.0351653806988
$ws.Range("C44").Value = 110.113382810351
$ws.Range("C45").Value = 110.3167133881492
$ws.Range("C46").Value = 109.4584115422884
$ws.Range("C47").Value = 109.7258314331384
$ws.Range("C48").Value = 109.9134098877834
$ws.Range("C49").Value = 110.3197463891681
$ws.Range("C50").Value = 110.6827680711938
$ws.Range("C51").Value = 111.4613489996468
$ws.Range("C52").Value = 112.0354998306003
$ws.Range("C53").Value = 113.0765161999174
$ws.Range("C54").Value = 113.3019597943043
$ws.Range("C55").Value = 114.3246246435287
$ws.Range("C56").Value = 114.7360872220881
$ws.Range("C57").Value = 115.1926143345042
$ws.Range("C58").Value = 115.0672168681821
$ws.Range("C59").Value = 115.3115288725653
$ws.Range("C60").Value = 115.2780050377524
$ws.Range("C61").Value = 115.8405155885507
$ws.Range("C62").Value = 116.6852303540751
$ws.Range("C63").Value = 117.5066334531552
$ws.Range("C64").Value = 118.4273243608429
$ws.Range("C65").Value = 119.0361152343294
$ws.Range("C66").Value = 120.5573651314358
$ws.Range("C67").Value = 121.4065790342464
$ws.Range("C68").Value = 120.1997202968823
$ws.Range("C69").Value = 120.2320000179186
$ws.Range("C70").Value = 119.4286429273413
$ws.Range("C71").Value = 119.9761507821712
$ws.Range("C72").Value = 120.9647688835819
$ws.Range("C73").Value = 121.8831692078534
$ws.Range("C74").Value = 122.7084576593252
$ws.Range("C75").Value = 123.519749170043
$ws.Range("C76").Value = 123.567313132901
$ws.Range("C77").Value = 124.8284044731425
$ws.Range("C78").Value = 125.3330465733894
$ws.Range("C79").Value = 126.1936909529623
$ws.Range("C80").Value = 126.1149877939172
